$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A1').Value2 = 'mx_state'
$ws.Range('B1').Value2 = 'mx_municipality'
$ws.Range('C1').Value2 = 'n_matriculas'
$ws.Range('D1').Value2 = 'pct_matriculas'
$ws.Range('B7').Value2 = 'Amatenango De La Frontera'
$ws.Range('B9').Value2 = 'Bejucal De Ocampo'
$ws.Range('B27').Value2 = 'Hidalgo Del Parral'
$ws.Range('A30').Value2 = 'Ciudad De México'
$ws.Range('B41').Value2 = 'Pánuco De Coronado'
$ws.Range('A43').Value2 = 'Estado De México'
$ws.Range('B44').Value2 = 'Almoloya De Alquisiras'
$ws.Range('B45').Value2 = 'Almoloya De Juárez'
$ws.Range('B54').Value2 = 'Ecatepec De Morelos'
$ws.Range('B57').Value2 = 'Ixtapan De La Sal'
$ws.Range('B58').Value2 = 'Ixtapan Del Oro'
$ws.Range('B63').Value2 = 'Naucalpan De Juárez'
$ws.Range('B66').Value2 = 'San Felipe Del Progreso'
$ws.Range('B67').Value2 = 'San Martín De Las Pirámides'
$ws.Range('B68').Value2 = 'San Simón De Guerero'
$ws.Range('B73').Value2 = 'Tenango Del Valle'
$ws.Range('B78').Value2 = 'Tlalnepantla De Baz'
$ws.Range('B81').Value2 = 'Valle De Bravo'
$ws.Range('B82').Value2 = 'Valle De Chalco Solidaridad'
$ws.Range('A88').Value2 = 'Guanajuato'
$ws.Range('B99').Value2 = 'Acapulco De Juárez'
$ws.Range('B100').Value2 = 'Ajuchitlán Del Progreso'
$ws.Range('B102').Value2 = 'Ayutla De Los Libres'
$ws.Range('B103').Value2 = 'Chilapa De Álvarez'
$ws.Range('B104').Value2 = 'Coahuayutla De José María Izazaga'
$ws.Range('B105').Value2 = 'Coyuca De Catalán'
$ws.Range('B109').Value2 = 'Huitzuco De Los Figueroa'
$ws.Range('B110').Value2 = 'Iguala De La Independencia'
$ws.Range('B111').Value2 = 'Zihuatanejo De Azueta'
$ws.Range('B117').Value2 = 'Taxco De Alarcón'
$ws.Range('B119').Value2 = 'Tepecoacuilco De Trujano'
$ws.Range('B125').Value2 = 'Tulancingo De Bravo'
$ws.Range('B127').Value2 = 'Ahualulco De Mercado'
$ws.Range('B131').Value2 = 'Encarnación De Díaz'
$ws.Range('B134').Value2 = 'Lagos De Moreno'
$ws.Range('B136').Value2 = 'Tizapán El Alto'
$ws.Range('B173').Value2 = 'Puente De Ixtla'
$ws.Range('B176').Value2 = 'Tlaltizapán De Zapata'
$ws.Range('B183').Value2 = 'Ixtlán Del Río'
$ws.Range('B190').Value2 = 'Oaxaca De Juárez'
$ws.Range('B192').Value2 = 'San Dionisio Del Mar'
$ws.Range('B210').Value2 = 'Tlacolula De Matamoros'
$ws.Range('B211').Value2 = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range('B216').Value2 = 'Ayotoxco De Guerero'
$ws.Range('B223').Value2 = 'Cuetzalan Del Progreso'
$ws.Range('B229').Value2 = 'Huehuetlán El Chico'
$ws.Range('B231').Value2 = 'Izúcar De Matamoros'
$ws.Range('B238').Value2 = 'Tecali De Herrera'
$ws.Range('B243').Value2 = 'Tlacotepec De Benito Juárez'
$ws.Range('B250').Value2 = 'Cadereyta De Montes'
$ws.Range('B253').Value2 = 'Pinal De Amoles'
$ws.Range('B256').Value2 = 'Ciudad Del Maíz'
$ws.Range('B268').Value2 = 'Jalpa De Méndez'
$ws.Range('B284').Value2 = 'Boca Del Río'
$ws.Range('B286').Value2 = 'Cazones De Herrera'
$ws.Range('B289').Value2 = 'Cosamaloapan De Carpio'
$ws.Range('B292').Value2 = 'Hueyapan De Ocampo'
$ws.Range('B293').Value2 = 'Ignacio De La Llave'
$ws.Range('B295').Value2 = 'Ixhuatlán De Madero'
$ws.Range('B298').Value2 = 'Martínez De La Torre'
$ws.Range('B304').Value2 = 'Poza Rica De Hidalgo'
$ws.Range('B307').Value2 = 'Soledad De Doblado'
# Remove the trailing footer rows (325-329) that are no longer part of the data
$ws.Range('A325:D329').EntireRow.Delete() | Out-Null
